# 0.0.13 - Change Upload System
# Rename the second x/y/z/value header block (I1:L1) to the "2" variants
# (x2, y2, z2, value2) so the sheet now distinguishes the first load-point
# columns (E1:H1 => x1,y1,z1,value1) from the second load-point columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "x2"
$ws.Range("J1").Value = "y2"
$ws.Range("K1").Value = "z2"
$ws.Range("L1").Value = "value2"

# Move/update the active selection to K11 (matches the saved view state).
$ws.Range("K11").Select() | Out-Null
